# AutoDM.xlsx update:
# - Insert a new "key" column (1..8) in front of the existing Child/Parent table.
# - Duplicate the existing "B -> D" relationship row.
# - Append two new "D -> E" relationship rows.
# - Highlight (red font) the two "B -> D" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "key"
$ws.Range("B1").Value = "Child"
$ws.Range("C1").Value = "Parent"

# Data rows: key, Child, Parent
$data = @(
  @(1, "A", "B"),
  @(2, "A", "B"),
  @(3, "A", "C"),
  @(4, "B", "D"),
  @(5, "B", "D"),
  @(6, "X", "Y"),
  @(7, "D", "E"),
  @(8, "D", "E")
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $r = $r + 1
}

# Highlight the two duplicated "B -> D" rows (now rows 5 and 6) in red font.
$ws.Rows("5:6").Font.Color = 255

# Match the final selection left by the author.
$ws.Range("C9").Select()
